$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.908.81'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.632.53'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'211.58"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Value = "'23.48"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('E9').Value = '  -0.65%  '
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '1.863.96'
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('D13').Value = '1.636.79'
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('E14').Value = '  -1.22%  '
$ws.Range('D15').Value = "'0.563"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('D16').Value = "'65.38"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').Value = '27.909.79'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').Value = "'229.27"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.58%  '
$ws.Range('E19').Value = '  +1.72%  '
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = "'4.34"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.76%  '
$ws.Range('E23').Value = '  -3.45%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = "'154.57"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.93%  '
$ws.Range('D26').Value = "'6.89"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = "'15.52"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  -0.52%  '
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('E32').Value = '  +1.11%  '
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('D34').Value = '1.393.22'
$ws.Range('E34').Value = '  -0.80%  '
$ws.Range('E35').Value = '  +0.42%  '
$ws.Range('E36').Value = '  +10.09%  '
$ws.Range('E37').Value = '  -0.90%  '
$ws.Range('E38').Value = '  +1.21%  '
$ws.Range('D39').Value = "'0.559"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('E40').Value = '  -3.18%  '
$ws.Range('E42').Value = '  -1.01%  '
$ws.Range('D43').Value = "'1.84"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.65%  '
$ws.Range('D44').Value = "'65.79"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.98%  '
$ws.Range('E45').Value = '  -1.98%  '
$ws.Range('D46').Value = '1.773.47'
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('E47').Value = '  -3.17%  '
$ws.Range('D48').Value = "'88.66"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('E49').Value = '  +1.51%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.0505"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'7.63"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.25%  '
